$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "button on CAN" column F additions
$ws.Range("F8").Value = 300
$ws.Range("F10").Value = "yes"
$ws.Range("F11").Value = 30
$ws.Range("F12").Value = "-"

# Row 11 ("Test time [min]") now holds numeric minute values across the board
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 30

# Row 12 ("Failure : CAN buffer overflow after[min]") now holds "-" across the board
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"

# Match the author's final selection
$ws.Range("F13").Select() | Out-Null
